# Delete the "lesion_volume" column (column D) entirely.
# This shifts columns E:H left to become D:G, matching the target diff:
#   - dimension changes from A1:H52 to A1:G52
#   - D1 header "lesion_volume" is removed, and the former E/F/G/H headers
#     (incongruent, congruent, congruent_log, incongruent_log) shift left
#     to become the new D/E/F/G headers
#   - every data row's column D value is removed and E/F/G/H shift left
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(4).Delete()
